$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (A Lag row values)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0.17"
$ws.Range("B2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "-0.01"
$ws.Range("B3").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "-0.09"
$ws.Range("B4").Style = "Normal"

# Column C
$ws.Range("C2").Value = "44.29***"
$ws.Range("C3").Value = "2.21***"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0.98"
$ws.Range("C4").Style = "Normal"

# Column D
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "-0.89"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "0.46***"
$ws.Range("D4").Value = "0.82*"
